# Rename the original (and only) sheet to "voltmetro"
$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws1.Name = "voltmetro"

# Add a new sheet "amperometro" right after voltmetro
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "amperometro"

# Header row
$ws2.Range("A1").Value = "R(OHM)"
$ws2.Range("B1").Value = "V"
$ws2.Range("C1").Value = "I(A)"

# Data rows
$data = @(
    @(1, 2.11, 1058),
    @(2, 2.53, 864),
    @(3, 2.60, 650),
    @(4, 2.64, 529),
    @(5, 2.66, 444),
    @(6, 2.68, 387),
    @(7, 2.70, 336),
    @(8, 2.70, 299),
    @(9, 2.71, 272),
    @(10, 2.72, 248)
)

$row = 2
foreach ($r in $data) {
    $ws2.Cells.Item($row, 1).Value = $r[0]
    $ws2.Cells.Item($row, 2).Value = $r[1]
    $ws2.Cells.Item($row, 3).Value = $r[2]
    $row++
}

# Apply the 2-decimal number format to the B column measurements
$ws2.Range("B2:B11").NumberFormat = "0.00"

# Make "amperometro" the active sheet/tab
$ws2.Activate()
